$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while preserving it as text (avoids Excel
# auto-converting numeric-looking strings like "0.999" into numbers,
# and keeps the cell back at its original/default style afterwards).
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$updates = @(
    @{ Cell = 'D2'; Value = '76.761.70' },
    @{ Cell = 'E2'; Value = '  +0.38%  ' },
    @{ Cell = 'D3'; Value = '3.145.23' },
    @{ Cell = 'E3'; Value = '  +6.28%  ' },
    @{ Cell = 'D4'; Value = '0.999' },
    @{ Cell = 'E4'; Value = '  -0.10%  ' },
    @{ Cell = 'D5'; Value = '201.88' },
    @{ Cell = 'E5'; Value = '  +1.74%  ' },
    @{ Cell = 'D6'; Value = '626.42' },
    @{ Cell = 'E6'; Value = '  +4.93%  ' },
    @{ Cell = 'D7'; Value = '0.999' },
    @{ Cell = 'E7'; Value = '  -0.03%  ' },
    @{ Cell = 'E8'; Value = '  +7.05%  ' },
    @{ Cell = 'D9'; Value = '0.560' },
    @{ Cell = 'E9'; Value = '  +1.35%  ' },
    @{ Cell = 'D10'; Value = '0.478' },
    @{ Cell = 'E10'; Value = '  +7.81%  ' },
    @{ Cell = 'E11'; Value = '  +0.47%  ' },
    @{ Cell = 'E12'; Value = '  +7.48%  ' },
    @{ Cell = 'D13'; Value = '3.711.15' },
    @{ Cell = 'E13'; Value = '  +5.97%  ' },
    @{ Cell = 'D14'; Value = '29.80' },
    @{ Cell = 'E14'; Value = '  +4.39%  ' },
    @{ Cell = 'E15'; Value = '  +6.28%  ' },
    @{ Cell = 'D16'; Value = '76.628.29' },
    @{ Cell = 'E16'; Value = '  +0.31%  ' },
    @{ Cell = 'D17'; Value = '3.139.34' },
    @{ Cell = 'E17'; Value = '  +5.93%  ' },
    @{ Cell = 'E18'; Value = '  +0.65%  ' },
    @{ Cell = 'D19'; Value = '9.29' },
    @{ Cell = 'E19'; Value = '  +6.04%  ' },
    @{ Cell = 'D20'; Value = '2.77' },
    @{ Cell = 'E20'; Value = '  +20.89%  ' },
    @{ Cell = 'D21'; Value = '400.25' },
    @{ Cell = 'E21'; Value = '  +5.91%  ' },
    @{ Cell = 'D22'; Value = '4.58' },
    @{ Cell = 'E22'; Value = '  +5.75%  ' },
    @{ Cell = 'E23'; Value = '  +2.23%  ' },
    @{ Cell = 'D24'; Value = '3.291.45' },
    @{ Cell = 'E24'; Value = '  +6.41%  ' },
    @{ Cell = 'E25'; Value = '  +7.56%  ' },
    @{ Cell = 'D26'; Value = '73.76' },
    @{ Cell = 'E26'; Value = '  +1.73%  ' },
    @{ Cell = 'D27'; Value = '10.38' },
    @{ Cell = 'E27'; Value = '  +6.58%  ' },
    @{ Cell = 'E28'; Value = '  +0.18%  ' },
    @{ Cell = 'E29'; Value = '  +4.31%  ' },
    @{ Cell = 'D30'; Value = '0.994' },
    @{ Cell = 'E30'; Value = '  -0.61%  ' },
    @{ Cell = 'D31'; Value = '8.51' },
    @{ Cell = 'E31'; Value = '  -0.95%  ' },
    @{ Cell = 'E32'; Value = '  +4.75%  ' },
    @{ Cell = 'D33'; Value = '523.99' },
    @{ Cell = 'E33'; Value = '  +5.49%  ' },
    @{ Cell = 'D34'; Value = '1.96' },
    @{ Cell = 'E34'; Value = '  +7.10%  ' },
    @{ Cell = 'D35'; Value = '0.134' },
    @{ Cell = 'E35'; Value = '  +20.55%  ' },
    @{ Cell = 'D36'; Value = '21.86' },
    @{ Cell = 'E36'; Value = '  +7.29%  ' },
    @{ Cell = 'D37'; Value = '0.999' },
    @{ Cell = 'E37'; Value = '  -0.08%  ' },
    @{ Cell = 'D38'; Value = '163.49' },
    @{ Cell = 'E38'; Value = '  -1.14%  ' },
    @{ Cell = 'D39'; Value = '196.38' },
    @{ Cell = 'E39'; Value = '  +8.73%  ' },
    @{ Cell = 'B40'; Value = 'PolygonEcosystemToken' },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol' },
    @{ Cell = 'D40'; Value = '0.387' },
    @{ Cell = 'E40'; Value = '  -1.60%  ' },
    @{ Cell = 'B41'; Value = 'WhiteBITCoin' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt' },
    @{ Cell = 'D41'; Value = '20.07' },
    @{ Cell = 'E41'; Value = '  +0.52%  ' },
    @{ Cell = 'E42'; Value = '  -5.37%  ' },
    @{ Cell = 'D43'; Value = '5.37' },
    @{ Cell = 'E43'; Value = '  +8.61%  ' },
    @{ Cell = 'E44'; Value = '  +0.06%  ' },
    @{ Cell = 'D45'; Value = '0.812' },
    @{ Cell = 'E45'; Value = '  +21.40%  ' },
    @{ Cell = 'E46'; Value = '  +8.34%  ' },
    @{ Cell = 'D47'; Value = '1.72' },
    @{ Cell = 'E47'; Value = '  +3.93%  ' },
    @{ Cell = 'D48'; Value = '41.90' },
    @{ Cell = 'E48'; Value = '  +5.00%  ' },
    @{ Cell = 'D49'; Value = '2.49' },
    @{ Cell = 'E49'; Value = '  +7.59%  ' },
    @{ Cell = 'D50'; Value = '0.617' },
    @{ Cell = 'E50'; Value = '  +4.34%  ' },
    @{ Cell = 'D51'; Value = '4.03' },
    @{ Cell = 'E51'; Value = '  +2.63%  ' }
)

foreach ($u in $updates) {
    Set-TextValue $ws.Range($u.Cell) $u.Value
}

Write-Host "Applied $($updates.Count) cell updates"
